$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xl enum constants used below (well-known Excel COM values):
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlContinuous = 1
$xlThin = 2
$xlMedium = -4138
$xlCenter = -4108

# --- start from a clean slate -------------------------------------------------
$ws.Cells.ClearFormats()

# Quantile index column (A) is not used in the new layout.
$ws.Columns.Item(1).Clear()

# Insert a new, blank spacer row above the header row.
$ws.Rows.Item(1).Insert()

# --- column widths ----------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 15.83203125
$ws.Columns.Item(2).ColumnWidth = 15.83203125
$ws.Range($ws.Cells.Item(1,3), $ws.Cells.Item(1,16384)).ColumnWidth = 15.83203125

# --- row heights --------------------------------------------------------------
$ws.Rows.Item("1:8").RowHeight = 27

# --- shared look: Cambria Math font, white fill, centred --------------------
$full = $ws.Range("B1:J8")
$full.Font.Name = "Cambria Math"
$full.Interior.Color = 16777215
$full.HorizontalAlignment = $xlCenter
$full.VerticalAlignment = $xlCenter

# --- bold emphasis: label column + header row --------------------------------
$ws.Range("B1:B8").Font.Bold = $true
$ws.Range("C2:J2").Font.Bold = $true

# --- borders per row ----------------------------------------------------------
# Row 1 (spacer row): bottom = medium
$row1 = $ws.Range("B1:J1")
$row1.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
$row1.Borders.Item($xlEdgeBottom).Weight = $xlMedium

# Row 2 (header row): top + bottom = medium
$row2 = $ws.Range("B2:J2")
$row2.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
$row2.Borders.Item($xlEdgeTop).Weight = $xlMedium
$row2.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
$row2.Borders.Item($xlEdgeBottom).Weight = $xlMedium

# Rows 4, 6, 7 (banded data rows): top + bottom = thin
foreach ($r in 4, 6, 7) {
    $row = $ws.Range("B$r`:J$r")
    $row.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
    $row.Borders.Item($xlEdgeTop).Weight = $xlThin
    $row.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
    $row.Borders.Item($xlEdgeBottom).Weight = $xlThin
}

# Row 8 (last row): bottom = medium
$row8 = $ws.Range("B8:J8")
$row8.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
$row8.Borders.Item($xlEdgeBottom).Weight = $xlMedium

# --- selection ------------------------------------------------------------
$ws.Range("C4").Select()

Write-Output "done"
